# Apply the cryptos-list refresh described by the commit.
# The workbook stores Price/Volume columns as plain text (inline strings),
# so for any new Price value that looks like a plain number we force the
# cell to Text format first -- this stops Excel from silently coercing
# values such as "1.00" or "0.110" into numbers (which would drop the
# formatting/trailing zeros the source data relies on).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "58.726.30"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "3.154.84"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "531.70"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "140.19"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.537"
$ws.Range("E8").Value = "  +16.23%  "
$ws.Range("D9").Value = "7.33"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "0.434"
$ws.Range("E10").Value = "  +5.82%  "
$ws.Range("D11").Value = "0.110"
$ws.Range("E11").Value = "  +2.39%  "
$ws.Range("E12").Value = "  +2.54%  "
$ws.Range("D13").Value = "3.705.19"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").Value = "26.18"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Value = "0.0000172"
$ws.Range("E15").Value = "  +4.47%  "
$ws.Range("D16").Value = "58.774.89"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "6.26"
$ws.Range("E17").Value = "  +4.16%  "
$ws.Range("D18").Value = "3.164.17"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").Value = "13.06"
$ws.Range("E19").Value = "  +2.95%  "
$ws.Range("D20").Value = "8.18"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").Value = "373.12"
$ws.Range("E21").Value = "  +4.83%  "
$ws.Range("D22").Value = "5.80"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "70.16"
$ws.Range("E24").Value = "  +1.72%  "
$ws.Range("D25").Value = "0.521"
$ws.Range("E25").Value = "  +3.08%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "0.995"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").Value = "8.24"
$ws.Range("E28").Value = "  +12.77%  "
$ws.Range("D29").Value = "0.0₃0869"
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "6.17"
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "22.21"
$ws.Range("E31").Value = "  +3.75%  "
$ws.Range("D32").Value = "1.88"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "5.20"
$ws.Range("E33").Value = "  +3.86%  "
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "159.29"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "6.28"
$ws.Range("E36").Value = "  +3.22%  "
$ws.Range("E37").Value = "  +5.89%  "
$ws.Range("D38").Value = "25.18"
$ws.Range("E38").Value = "  -2.92%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "2.655.57"
$ws.Range("E39").Value = "  +10.17%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.68"
$ws.Range("E40").Value = "  -2.32%  "
$ws.Range("E41").Value = "  +2.32%  "
$ws.Range("D42").Value = "4.24"
$ws.Range("E42").Value = "  +5.55%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "38.81"
$ws.Range("E43").Value = "  +2.79%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "0.712"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0286"
$ws.Range("E45").Value = "  +6.66%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "3.198.60"
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("E48").Value = "  +13.86%  "
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("D50").Value = "6.20"
$ws.Range("E50").Value = "  +2.60%  "
$ws.Range("D51").Value = "20.35"
$ws.Range("E51").Value = "  +2.13%  "
